# Appends the next batch of sensor log rows (2026-01-28, ~12:24-12:25)
# captured by the logger to the PIR, Humidity and Temperature sheets.
#
# Note: the Date column values are entered with a leading "'" so Excel
# stores them as literal text ("2026-01-28") instead of auto-converting
# them to a real date serial, matching how the rest of the log is stored.
# The Humidity sheet's percentage values ("86.4%" etc.) need the same
# treatment so they stay literal text instead of becoming numeric percents.
$wb = $excel.ActiveWorkbook
$dateStr = "2026-01-28"

$ws = $wb.Worksheets.Item("PIR")
$rows = @(
    @("12:24:24", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("12:24:29", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("12:24:35", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("12:24:39", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("12:24:44", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("12:24:49", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("12:24:55", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("12:24:59", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("12:25:04", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("12:25:09", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("12:25:16", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("12:25:19", "12:00", "Bathroom", "No Motion", "Inactive")
)
$startRow = 240
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = "'" + $dateStr
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
}

$ws = $wb.Worksheets.Item("Humidity")
$rows = @(
    @("12:24:26", "12:00", "Bathroom", "86.4%", "Active"),
    @("12:24:30", "12:00", "Bathroom", "87.3%", "Active"),
    @("12:24:34", "12:00", "Bathroom", "87.3%", "Active"),
    @("12:24:38", "12:00", "Bathroom", "86.4%", "Active"),
    @("12:24:42", "12:00", "Bathroom", "87.3%", "Active"),
    @("12:24:46", "12:00", "Bathroom", "86.4%", "Active"),
    @("12:24:54", "12:00", "Bathroom", "87.3%", "Active"),
    @("12:24:58", "12:00", "Bathroom", "86.3%", "Active"),
    @("12:25:06", "12:00", "Bathroom", "86.4%", "Active"),
    @("12:25:10", "12:00", "Bathroom", "87.2%", "Active"),
    @("12:25:14", "12:00", "Bathroom", "87.3%", "Active"),
    @("12:25:18", "12:00", "Bathroom", "86.4%", "Active")
)
$startRow = 226
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = "'" + $dateStr
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = "'" + $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
}

$ws = $wb.Worksheets.Item("Temperature")
$rows = @(
    @("12:24:26", "12:00", "Bathroom", "23.0C", "Active"),
    @("12:24:30", "12:00", "Bathroom", "23.0C", "Active"),
    @("12:24:35", "12:00", "Bathroom", "23.0C", "Active"),
    @("12:24:38", "12:00", "Bathroom", "23.0C", "Active"),
    @("12:24:42", "12:00", "Bathroom", "23.0C", "Active"),
    @("12:24:46", "12:00", "Bathroom", "23.0C", "Active"),
    @("12:24:55", "12:00", "Bathroom", "23.0C", "Active"),
    @("12:24:59", "12:00", "Bathroom", "23.0C", "Active"),
    @("12:25:07", "12:00", "Bathroom", "23.0C", "Active"),
    @("12:25:11", "12:00", "Bathroom", "22.9C", "Active"),
    @("12:25:15", "12:00", "Bathroom", "23.0C", "Active"),
    @("12:25:19", "12:00", "Bathroom", "23.0C", "Active")
)
$startRow = 226
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = "'" + $dateStr
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
}
